# Updating the Forecast Portfolio
# Shift all timestamps in column A (rows 2-97) forward by 16 days,
# and set the new solar production forecast values in column B (rows 30-40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all the timestamps (column A, rows 2 through 97) forward by 16 days.
for ($row = 2; $row -le 97; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value2 = $cell.Value2 + 16
}

# Update the forecasted solar production values (column B, rows 30 through 40).
$newValues = @{
    30 = 5
    31 = 39
    32 = 98
    33 = 176
    34 = 288
    35 = 395
    36 = 538
    37 = 674
    38 = 784
    39 = 874
    40 = 957
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}
